$d = $word.ActiveDocument

# Delete paragraphs 1 (Title) and 2 (blank) by deleting the range spanning both.
$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs(2).Range.End
$r = $d.Range($start, $end)
$r.Delete()

$count = $d.Paragraphs.Count
Write-Output "count=$count"
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    Write-Output "[$i] len=$($t.Length) :: $t"
}
